$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.425.81'
$ws.Range("E2").Value = '  -2.33%  '
$ws.Range("D3").Value = '3.409.91'
$ws.Range("E3").Value = '  -3.38%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '''579.04'
$ws.Range("E5").Value = '  -4.11%  '
$ws.Range("D6").Value = '''133.66'
$ws.Range("E6").Value = '  -6.96%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '3.408.21'
$ws.Range("E8").Value = '  -3.42%  '
$ws.Range("D9").Value = '''0.479'
$ws.Range("E9").Value = '  -5.94%  '
$ws.Range("E10").Value = '  -9.00%  '
$ws.Range("D11").Value = '''6.99'
$ws.Range("E11").Value = '  -10.17%  '
$ws.Range("D12").Value = '''0.371'
$ws.Range("E12").Value = '  -8.32%  '
$ws.Range("D13").Value = '3.995.98'
$ws.Range("E13").Value = '  -3.25%  '
$ws.Range("D14").Value = '''0.0000176'
$ws.Range("E14").Value = '  -8.88%  '
$ws.Range("D15").Value = '3.419.63'
$ws.Range("E15").Value = '  -2.50%  '
$ws.Range("E16").Value = '  -1.77%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '64.469.43'
$ws.Range("E17").Value = '  -2.14%  '
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").Value = '''25.78'
$ws.Range("E18").Value = '  -8.78%  '
$ws.Range("D19").Value = '''9.35'
$ws.Range("E19").Value = '  -14.51%  '
$ws.Range("D20").Value = '''5.69'
$ws.Range("E20").Value = '  -7.79%  '
$ws.Range("D21").Value = '''13.35'
$ws.Range("E21").Value = '  -8.26%  '
$ws.Range("D22").Value = '''376.63'
$ws.Range("E22").Value = '  -10.34%  '
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").Value = '''0.537'
$ws.Range("E24").Value = '  -9.51%  '
$ws.Range("E25").Value = '  -7.23%  '
$ws.Range("D26").Value = '3.548.22'
$ws.Range("E26").Value = '  -3.35%  '
$ws.Range("D27").Value = '''0.0000103'
$ws.Range("E27").Value = '  -8.90%  '
$ws.Range("D29").Value = '''7.11'
$ws.Range("E29").Value = '  -8.93%  '
$ws.Range("D30").Value = '''2.16'
$ws.Range("E30").Value = '  -12.15%  '
$ws.Range("D31").Value = '''7.88'
$ws.Range("E31").Value = '  -10.98%  '
$ws.Range("D32").Value = '3.431.47'
$ws.Range("E32").Value = '  -2.98%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").Value = '''22.79'
$ws.Range("E34").Value = '  -5.69%  '
$ws.Range("D35").Value = '''0.139'
$ws.Range("E35").Value = '  -9.79%  '
$ws.Range("D36").Value = '''168.32'
$ws.Range("E36").Value = '  -4.58%  '
$ws.Range("D37").Value = '''1.16'
$ws.Range("E37").Value = '  -14.17%  '
$ws.Range("D38").Value = '''6.60'
$ws.Range("E38").Value = '  -12.05%  '
$ws.Range("D39").Value = '''1.42'
$ws.Range("E39").Value = '  -11.18%  '
$ws.Range("D40").Value = '''4.52'
$ws.Range("E40").Value = '  -13.49%  '
$ws.Range("D41").Value = '''0.0747'
$ws.Range("E41").Value = '  -8.05%  '
$ws.Range("D42").Value = '''0.802'
$ws.Range("E42").Value = '  -6.29%  '
$ws.Range("E43").Value = '  +0.20%  '
$ws.Range("D44").Value = '''41.39'
$ws.Range("E44").Value = '  -8.92%  '
$ws.Range("D45").Value = '''4.22'
$ws.Range("E45").Value = '  -14.63%  '
$ws.Range("E46").Value = '  -10.27%  '
$ws.Range("E47").Value = '  -1.65%  '
$ws.Range("D48").Value = '''22.17'
$ws.Range("E48").Value = '  -4.78%  '
$ws.Range("D49").Value = '''6.40'
$ws.Range("E49").Value = '  -8.74%  '
$ws.Range("D50").Value = '2.170.12'
$ws.Range("E50").Value = '  -6.48%  '
$ws.Range("D51").Value = '''1.94'
$ws.Range("E51").Value = '  -17.73%  '
